# Add files via upload
#
# The new "image_base64" column header is inserted at J1, pushing the
# existing image_base64_1 .. image_base64_4 headers one column to the
# right (image_base64_5 drops off the end, past N1 which is the last
# used column).
#
# Shift right-to-left so we never clobber a value before reading it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = $ws.Range("M1").Value2
$ws.Range("M1").Value = $ws.Range("L1").Value2
$ws.Range("L1").Value = $ws.Range("K1").Value2
$ws.Range("K1").Value = $ws.Range("J1").Value2
$ws.Range("J1").Value = "image_base64"

$ws.Range("N1").Select()
